$d = $word.ActiveDocument
$section = $d.Sections.First
$header = $section.Headers.Item(2)  # wdHeaderFooterFirstPage = 2
$header.Range.Find.Execute("Dr. med. Thiên-Trí Lâm", $true, $false, $false, $false, $false, $true, 1, $false, "PD Dr. med. Thiên-Trí Lâm", 2)
